# Weekly update: a new price record is inserted at the top of the
# Ciboulette price history (row 197) and every subsequent record shifts
# down by one row, with the last existing record (old row 285) now
# landing in a brand-new row 286.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 197
$lastRow  = 285
$newLastRow = $lastRow + 1

# --- 1) snapshot the "old" values for the columns that shift (D, J, K, L, M, P)
#        before we start overwriting anything.
$dVals = @{}
$jVals = @{}
$kVals = @{}
$lVals = @{}
$mVals = @{}
$pVals = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dVals[$r] = $ws.Cells.Item($r, 4).Value2()
    $jVals[$r] = $ws.Cells.Item($r, 10).Value2()
    $kVals[$r] = $ws.Cells.Item($r, 11).Value2()
    $lVals[$r] = $ws.Cells.Item($r, 12).Value2()
    $mVals[$r] = $ws.Cells.Item($r, 13).Value2()
    $pVals[$r] = $ws.Cells.Item($r, 16).Value2()
}

$dateFormat = $ws.Cells.Item($firstRow, 4).NumberFormat

# --- 2) the newest weekly record goes into the first row; everything that
#        used to be there moves one row down.
$ws.Cells.Item($firstRow, 4).Value = 44636
$ws.Cells.Item($firstRow, 4).NumberFormat = $dateFormat
$ws.Cells.Item($firstRow, 10).Value = 160

# --- 3) shift the rest of the series down by one row (row r gets what used
#        to be in row r-1), including the brand-new last row.
for ($r = ($firstRow + 1); $r -le $newLastRow; $r++) {
    $src = $r - 1

    $ws.Cells.Item($r, 4).Value = $dVals[$src]
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    $ws.Cells.Item($r, 10).Value = $jVals[$src]
    $ws.Cells.Item($r, 11).Value = $kVals[$src]
    $ws.Cells.Item($r, 12).Value = $lVals[$src]
    $ws.Cells.Item($r, 13).Value = $mVals[$src]
    $ws.Cells.Item($r, 16).Value = $pVals[$src]
}

# --- 4) the new last row also needs the rest of the (constant) row data,
#        copied straight across from the row above it.
$ws.Cells.Item($newLastRow, 1).Value  = $ws.Cells.Item($newLastRow - 1, 1).Value()
$ws.Cells.Item($newLastRow, 2).Value  = $ws.Cells.Item($newLastRow - 1, 2).Value()
$ws.Cells.Item($newLastRow, 3).Value  = $ws.Cells.Item($newLastRow - 1, 3).Value()
$ws.Cells.Item($newLastRow, 5).Value  = $ws.Cells.Item($newLastRow - 1, 5).Value()
$ws.Cells.Item($newLastRow, 6).Value  = $ws.Cells.Item($newLastRow - 1, 6).Value()
$ws.Cells.Item($newLastRow, 7).Value  = $ws.Cells.Item($newLastRow - 1, 7).Value()
$ws.Cells.Item($newLastRow, 8).Value  = $ws.Cells.Item($newLastRow - 1, 8).Value()
$ws.Cells.Item($newLastRow, 9).Value  = $ws.Cells.Item($newLastRow - 1, 9).Value()
$ws.Cells.Item($newLastRow, 14).Value = $ws.Cells.Item($newLastRow - 1, 14).Value()
$ws.Cells.Item($newLastRow, 15).Value = $ws.Cells.Item($newLastRow - 1, 15).Value()
$ws.Cells.Item($newLastRow, 17).Value = $ws.Cells.Item($newLastRow - 1, 17).Value()
$ws.Cells.Item($newLastRow, 18).Value = $ws.Cells.Item($newLastRow - 1, 18).Value()
